# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N ("Late") on the
# "Repayment Schedule" sheet, widen it to match the author's manual
# resize, and leave the new active selection/tab on that sheet - matching
# the author's workflow of adding an extra instalment-variance column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column at position N (14th column). Excel shifts the
# existing N:P data (Late / heading / Outstanding) right to O:Q and
# leaves the freshly inserted column N empty.
$ws.Columns.Item(14).Insert()

# The author then manually widened the new blank column to 10 characters.
$ws.Columns.Item(14).ColumnWidth = 9.166666666666666

# Reflect the author's final selection/active sheet: they ended up with
# the cursor on S8 of the Repayment Schedule tab (making it the active
# tab of the workbook).
$ws.Range("S8").Select()
